# Agregar planes de Residencial de 50mb (Plans sheet)
$wb = $excel.ActiveWorkbook

$plans = $wb.Worksheets.Item("Plans")

# Row 5 -> Residencial / Con_TotalPlay_TV / 50 / Amazon Prime = Si
$plans.Range("B5").Value = "Residencial"
$plans.Range("C5").Value = "Con_TotalPlay_TV"
$plans.Range("D5").Value = 50
$plans.Range("E5").Value = "Si"

# Row 6 -> Residencial / Con_TotalPlay_TV / 50 / Netflix = Netflix Estandar
$plans.Range("B6").Value = "Residencial"
$plans.Range("C6").Value = "Con_TotalPlay_TV"
$plans.Range("D6").Value = 50
$plans.Range("F6").Value = "Netflix Estándar"

# Row 7 -> Residencial / Sin_TotalPlay_TV / 50 / Amazon Prime = Si
$plans.Range("B7").Value = "Residencial"
$plans.Range("C7").Value = "Sin_TotalPlay_TV"
$plans.Range("D7").Value = 50
$plans.Range("E7").Value = "Si"

# Row 8 -> Residencial / Sin_TotalPlay_TV / 50 / Netflix = Netflix Premium
$plans.Range("B8").Value = "Residencial"
$plans.Range("C8").Value = "Sin_TotalPlay_TV"
$plans.Range("D8").Value = 50
$plans.Range("F8").Value = "Netflix Premium"

# Row 9 -> cleared (previously had the 5th Micronegocio plan)
$plans.Range("B9").Value = ""
$plans.Range("C9").Value = ""
$plans.Range("D9").Value = ""

$plans.Range("E10").Select()
